$wb = $excel.ActiveWorkbook

$wsPlan = $wb.Worksheets.Item("plan")
$wsExecution = $wb.Worksheets.Item("execution")

# Sprint 2 planning: the planning start date moved from 10-26-2018 to 10-18-2018.
$wsPlan.Range("E1").Value = "1: 10-18-2018"
$wsExecution.Range("E1").Value = "1: 10-18-2018"

# Remove the stale chart-tracking defined names left over from the old burndown chart.
while ($wb.Names.Count -gt 0) {
    $wb.Names.Item(1).Delete()
}

# Restore the plan sheet's selection to a single cell, then move the active
# selection/tab onto the execution sheet for sprint 2 planning.
$wsPlan.Range("E1").Select()
$wsExecution.Range("E1").Select()
